$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1810.25
$ws.Range("I28").Value = 1060.5625
$ws.Range("J28").Value = 3309.625
$ws.Range("K28").Value = 1060.5625
$ws.Range("L28").Value = 3309.625
$ws.Range("M28").Value = -575.5625
$ws.Range("N28").Value = -4279.625

# Row 64
$ws.Range("H64").Value = 3132.8
$ws.Range("I64").Value = 2917
$ws.Range("J64").Value = 3276.6667
$ws.Range("K64").Value = 2917
$ws.Range("L64").Value = 3276.6667
$ws.Range("M64").Value = -2669
$ws.Range("N64").Value = -3772.6667

# Row 67
$ws.Range("H67").Value = 3132.8
$ws.Range("I67").Value = 2917
$ws.Range("J67").Value = 3276.6667
$ws.Range("K67").Value = 2917
$ws.Range("L67").Value = 3276.6667
$ws.Range("M67").Value = -2059
$ws.Range("N67").Value = -4992.6667

# Row 76
$ws.Range("H76").Value = 2982.2942
$ws.Range("I76").Value = 2976.8462
$ws.Range("K76").Value = 2976.8462
$ws.Range("M76").Value = -2661.8462

# Row 79
$ws.Range("H79").Value = 2982.2942
$ws.Range("I79").Value = 2976.8462
$ws.Range("K79").Value = 2976.8462
$ws.Range("M79").Value = -1884.8462

# Row 103
$ws.Range("H103").Value = 3427.8333
$ws.Range("I103").Value = 17368
$ws.Range("J103").Value = 639.8
$ws.Range("K103").Value = 52104
$ws.Range("L103").Value = 1919.4
$ws.Range("M103").Value = -51518
$ws.Range("N103").Value = -3091.4

# Row 121
$ws.Range("H121").Value = 3821.2
$ws.Range("J121").Value = 4701.25
$ws.Range("L121").Value = 14103.75
$ws.Range("N121").Value = -17597.75

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7863.617
$ws.Range("I32").Value = 5020.294
$ws.Range("K32").Value = 5020.294
$ws.Range("M32").Value = -4733.294

# Row 110
$ws.Range("H110").Value = 2351
$ws.Range("I110").Value = 2337.2632
$ws.Range("J110").Value = 2403.2
$ws.Range("K110").Value = 2337.2632
$ws.Range("L110").Value = 2403.2
$ws.Range("M110").Value = -292.2631999999999
$ws.Range("N110").Value = -6493.2

# Row 122
$ws.Range("H122").Value = 2465
$ws.Range("I122").Value = 2412.5
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 7237.5
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4787.5
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3262.4167
$ws.Range("I86").Value = 3533.3333
$ws.Range("J86").Value = 3172.111
$ws.Range("K86").Value = 3533.3333
$ws.Range("L86").Value = 3172.111
$ws.Range("M86").Value = -2410.3333
$ws.Range("N86").Value = -5418.111

# Row 89
$ws.Range("H89").Value = 3262.4167
$ws.Range("I89").Value = 3533.3333
$ws.Range("J89").Value = 3172.111
$ws.Range("K89").Value = 17666.6665
$ws.Range("L89").Value = 15860.555
$ws.Range("M89").Value = -12050.6665
$ws.Range("N89").Value = -27092.555

# Row 105
$ws.Range("H105").Value = 2515.3845
$ws.Range("I105").Value = 2516.6667
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2516.6667
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -769.6667000000002
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2030.8518
$ws.Range("I31").Value = 1251.75
$ws.Range("J31").Value = 2654.1333
$ws.Range("K31").Value = 1251.75
$ws.Range("L31").Value = 2654.1333
$ws.Range("M31").Value = -956.75
$ws.Range("N31").Value = -3244.1333

# Row 34
$ws.Range("H34").Value = 2030.8518
$ws.Range("I34").Value = 1251.75
$ws.Range("J34").Value = 2654.1333
$ws.Range("K34").Value = 1251.75
$ws.Range("L34").Value = 2654.1333
$ws.Range("M34").Value = -1049.75
$ws.Range("N34").Value = -3058.1333

# Row 62
$ws.Range("H62").Value = 2476.6316
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 2672.889
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 2672.889
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -3920.889

# Row 65
$ws.Range("H65").Value = 2476.6316
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 2672.889
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 13364.445
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -19604.445

# Row 141
$ws.Range("H141").Value = 51304.715
$ws.Range("J141").Value = 51304.715
$ws.Range("L141").Value = 51304.715
$ws.Range("N141").Value = -61664.715

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 224.64285
$ws.Range("J12").Value = 189.88235
$ws.Range("L12").Value = 569.64705
$ws.Range("N12").Value = -915.64705

# Row 98
$ws.Range("H98").Value = 1747.7142
$ws.Range("J98").Value = 1500
$ws.Range("L98").Value = 4500
$ws.Range("N98").Value = -7496

# Row 131
$ws.Range("H131").Value = 887.53
$ws.Range("J131").Value = 902.46313
$ws.Range("L131").Value = 2707.38939
$ws.Range("N131").Value = -12787.38939

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 112556.336
$ws.Range("I122").Value = 112556.336
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 337669.008
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -335219.008
$ws.Range("N122").ClearContents()

# Row 123
$ws.Range("H123").Value = 15223.692
$ws.Range("J123").Value = 15223.692
$ws.Range("L123").Value = 15223.692
$ws.Range("N123").Value = -20123.692

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 91811
$ws.Range("I22").Value = 200636.2
$ws.Range("J22").Value = 1123.3334
$ws.Range("K22").Value = 200636.2
$ws.Range("L22").Value = 1123.3334
$ws.Range("M22").Value = -200341.2
$ws.Range("N22").Value = -1713.3334

# Row 27
$ws.Range("H27").Value = 91811
$ws.Range("I27").Value = 200636.2
$ws.Range("J27").Value = 1123.3334
$ws.Range("K27").Value = 200636.2
$ws.Range("L27").Value = 1123.3334
$ws.Range("M27").Value = -200529.2
$ws.Range("N27").Value = -1337.3334

# Row 122
$ws.Range("H122").Value = 4993.143
$ws.Range("I122").Value = 5337.375
$ws.Range("J122").Value = 4534.1665
$ws.Range("K122").Value = 16012.125
$ws.Range("L122").Value = 13602.4995
$ws.Range("M122").Value = -13562.125
$ws.Range("N122").Value = -18502.4995

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 64932.473
$ws.Range("I122").Value = 1515.2858
$ws.Range("J122").Value = 101925.836
$ws.Range("K122").Value = 4545.857400000001
$ws.Range("L122").Value = 305777.508
$ws.Range("M122").Value = -2095.857400000001
$ws.Range("N122").Value = -310677.508

# Row 135
$ws.Range("H135").Value = 35000
$ws.Range("J135").Value = 35000
$ws.Range("L135").Value = 35000
$ws.Range("N135").Value = -45140
